{"js": "// Remove the whole list paragraph \"Finish Desert and North Pole\"\n// (the task soon got split into other items, per the commit message).\nconst body = context.document.body;\nconst results = body.search(\"Finish Desert and North Pole\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const para = results.items[i].paragraphs.getFirst();\n  para.delete();\n}\nawait context.sync();\n", "ps1": "# Remove the whole list-paragraph \"Finish Desert and North Pole\" \u2014 the\n# task was superseded/merged into other items (see commit message).\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"Finish Desert and North Pole\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
